# Clear particles at the start of optimization load.
# Update the "studio" sheet's saved UI state so that, when the project is
# reopened, it comes up on the "optimize" tool tab with the "Groomed" view
# (instead of stale "analysis" / "Reconstructed" state) and a refreshed
# zoom_state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studio")

$ws.Range("B5").Value = "optimize"
$ws.Range("B6").Value = "Groomed"

# zoom_state is stored as text (not a number) in the workbook, so force a
# text entry the same way typing '2 into the cell in Excel would.
$ws.Range("B7").Value = "'2"
